$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.85
$ws.Range("I3").Value = 4.33
$ws.Range("J3").Value = 2.6
$ws.Range("M3").Value = 1.08
$ws.Range("N3").Value = 7.5
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62
$ws.Range("AC3").Value = 7.5
$ws.Range("AG3").Value = 1250
$ws.Range("AQ3").Value = 41
$ws.Range("AR3").Value = 67
$ws.Range("BB3").Value = 351
